$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '36.712.10'
$ws.Range("E2").Value = '  -1.11%  '

$ws.Range("D3").Value = "'" + '2.081.44'
$ws.Range("E3").Value = '  +1.44%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").Value = "'" + '244.66'
$ws.Range("E5").Value = '  -1.48%  '

$ws.Range("D6").Value = "'" + '0.649'
$ws.Range("E6").Value = '  -2.17%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("E8").Value = '  -6.96%  '

$ws.Range("D9").Value = "'" + '58.73'
$ws.Range("E9").Value = '  -1.94%  '

$ws.Range("D10").Value = "'" + '0.364'
$ws.Range("E10").Value = '  -4.57%  '

$ws.Range("D11").Value = "'" + '0.0759'
$ws.Range("E11").Value = '  -2.21%  '

$ws.Range("E12").Value = '  +1.01%  '

$ws.Range("D13").Value = "'" + '14.89'
$ws.Range("E13").Value = '  -6.76%  '

$ws.Range("E14").Value = '  +0.75%  '

$ws.Range("D15").Value = "'" + '2.386.60'
$ws.Range("E15").Value = '  +1.59%  '

$ws.Range("D16").Value = "'" + '5.46'
$ws.Range("E16").Value = '  -4.27%  '

$ws.Range("D17").Value = "'" + '2.107.81'
$ws.Range("E17").Value = '  +2.74%  '

$ws.Range("D18").Value = "'" + '36.682.73'
$ws.Range("E18").Value = '  -1.27%  '

$ws.Range("D19").Value = "'" + '17.16'
$ws.Range("E19").Value = '  -4.49%  '

$ws.Range("D20").Value = "'" + '72.36'
$ws.Range("E20").Value = '  -3.37%  '

$ws.Range("D21").Value = "'" + '0.0₃0874'
$ws.Range("E21").Value = '  -1.95%  '

$ws.Range("D22").Value = "'" + '5.41'
$ws.Range("E22").Value = '  +0.73%  '

$ws.Range("D23").Value = "'" + '239.65'
$ws.Range("E23").Value = '  +1.00%  '

$ws.Range("E24").Value = '  +0.03%  '

$ws.Range("D25").Value = "'" + '2.39'
$ws.Range("E25").Value = '  -3.37%  '

$ws.Range("D26").Value = "'" + '9.79'
$ws.Range("E26").Value = '  +2.96%  '

$ws.Range("E27").Value = '  -1.30%  '

$ws.Range("D28").Value = "'" + '167.03'
$ws.Range("E28").Value = '  -1.36%  '

$ws.Range("D29").Value = "'" + '20.58'
$ws.Range("E29").Value = '  +2.63%  '

$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = "'" + '0.123'
$ws.Range("E30").Value = '  -1.16%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = "'" + '5.29'
$ws.Range("E31").Value = '  +9.84%  '

$ws.Range("D32").Value = "'" + '1.17'
$ws.Range("E32").Value = '  +4.33%  '

$ws.Range("E33").Value = '  +3.41%  '

$ws.Range("D34").Value = "'" + '0.0604'
$ws.Range("E34").Value = '  -2.37%  '

$ws.Range("E35").Value = '  +5.45%  '

$ws.Range("E36").Value = '  +0.18%  '

$ws.Range("E37").Value = '  +4.32%  '

$ws.Range("D38").Value = "'" + '0.0823'
$ws.Range("E38").Value = '  -8.23%  '

$ws.Range("D39").Value = "'" + '1.26'
$ws.Range("E39").Value = '  -5.94%  '

$ws.Range("D40").Value = "'" + '0.0219'
$ws.Range("E40").Value = '  -1.39%  '

$ws.Range("E41").Value = '  +0.74%  '

$ws.Range("D42").Value = "'" + '4.83'
$ws.Range("E42").Value = '  -6.37%  '

$ws.Range("E43").Value = '  +0.94%  '

$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").Value = "'" + '2.85'
$ws.Range("E44").Value = '  -11.13%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = "'" + '95.85'
$ws.Range("E45").Value = '  -0.18%  '

$ws.Range("D46").Value = "'" + '15.95'
$ws.Range("E46").Value = '  -7.92%  '

$ws.Range("D47").Value = "'" + '1.362.92'
$ws.Range("E47").Value = '  +6.90%  '

$ws.Range("D48").Value = "'" + '7.29'
$ws.Range("E48").Value = '  +6.93%  '

$ws.Range("D49").Value = "'" + '2.42'
$ws.Range("E49").Value = '  -0.54%  '

$ws.Range("E50").Value = '  +0.95%  '

$ws.Range("D51").Value = "'" + '2.267.43'
$ws.Range("E51").Value = '  +1.53%  '
